$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New 5th column header: "字符串相加" (leetcode "Add Strings"), mirroring
# the existing D1 "相交链表" header (hyperlinked, styled with the built-in
# Hyperlink cell style).
$ws.Range("E1").Value = "字符串相加"

# Adding via the Hyperlinks collection both creates the relationship and
# (as a side effect) stamps E1 with the hyperlink style - matches what
# happened to D1 originally.
$ws.Hyperlinks.Add($ws.Range("E1"), "https://leetcode-cn.com/problems/add-strings/solution/zi-fu-chuan-xiang-jia-by-leetcode-solution/")

# Re-apply D1's exact formatting (same cell style xf) onto E1 so both
# header cells share one style entry instead of Hyperlinks.Add minting a
# fresh one.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New column width to fit the header text.
$ws.Columns.Item(5).ColumnWidth = 10.2857142857142

# Selection moves on to F5 after the edit.
$ws.Range("F5").Select() | Out-Null
